# formDef.json choice lists cannot contain integers in the "value" column,
# so every numeric data_value in the "choices" sheet (column B, rows 2-95)
# is rewritten as text, prefixed with "a" (e.g. 1 -> "a1", -777 -> "a-777").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("choices")

$valueMap = @{2="a1"; 3="a0"; 4="a888"; 5="a1"; 6="a0"; 7="a888"; 8="a888"; 9="a999"; 10="a1"; 11="a2"; 12="a1"; 13="a2"; 14="a3"; 15="a4"; 16="a5"; 17="a6"; 18="a7"; 19="a888"; 20="a1"; 21="a2"; 22="a3"; 23="a4"; 24="a5"; 25="a6"; 26="a7"; 27="a8"; 28="a9"; 29="a10"; 30="a888"; 31="a0"; 32="a1"; 33="a2"; 34="a3"; 35="a4"; 36="a5"; 37="a6"; 38="a7"; 39="a888"; 40="a0"; 41="a1"; 42="a2"; 43="a3"; 44="a4"; 45="a888"; 46="a0"; 47="a1"; 48="a2"; 49="a3"; 50="a4"; 51="a5"; 52="a888"; 53="a1"; 54="a2"; 55="a3"; 56="a4"; 57="a1"; 58="a2"; 59="a3"; 60="a999"; 61="a0"; 62="a1"; 63="a888"; 64="a999"; 65="a1"; 66="a0"; 67="a888"; 68="a999"; 69="a-777"; 70="a-888"; 71="a-999"; 72="a1"; 73="a0"; 74="a2"; 75="a999"; 76="a1"; 77="a2"; 78="a3"; 79="a1"; 80="a2"; 81="a3"; 82="a4"; 83="a5"; 84="a6"; 85="a888"; 86="a0"; 87="a1"; 88="a1"; 89="a1"; 90="a0"; 91="a1"; 92="a0"; 93="a1"; 94="a0"; 95="a2"}

foreach ($r in $valueMap.Keys) {
    $ws.Cells.Item($r, 2).Value = $valueMap[$r]
}

# B18/B19 were the only two cells using the right-aligned number style; now
# that they hold text like the rest of column B, drop the alignment so they
# pick up the same (left/general) style as every other data cell in column B.
$ws.Range("B18").HorizontalAlignment = 1
$ws.Range("B19").HorizontalAlignment = 1

# The user ended the session on the "choices" sheet, cell E6 selected.
$ws.Activate()
$ws.Range("E6").Select()
